$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.166.37'
$ws.Range("D3").Value = '1.677.39'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  -0.55%  '
$ws.Range("D5").Value = '''209.67'
$ws.Range("E5").Value = '  -4.35%  '
$ws.Range("D6").Value = '''0.5280'
$ws.Range("E6").Value = '  -4.50%  '
$ws.Range("D7").Value = '''1.005'
$ws.Range("E7").Value = '  -0.57%  '
$ws.Range("D8").Value = '''0.2675'
$ws.Range("E8").Value = '  -1.64%  '
$ws.Range("D9").Value = '''0.06279'
$ws.Range("E9").Value = '  -3.21%  '
$ws.Range("D10").Value = '''21.20'
$ws.Range("E10").Value = '  -4.28%  '
$ws.Range("D11").Value = '''0.07519'
$ws.Range("E11").Value = '  -1.34%  '
$ws.Range("D12").Value = '1.703.27'
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").Value = '''4.476'
$ws.Range("E13").Value = '  -1.93%  '
$ws.Range("D14").Value = '''0.5636'
$ws.Range("E14").Value = '  -3.29%  '
$ws.Range("D15").Value = '''0.000008092'
$ws.Range("E15").Value = '  -4.48%  '
$ws.Range("D16").Value = '''66.13'
$ws.Range("E16").Value = '  +1.16%  '
$ws.Range("D17").Value = '26.168.25'
$ws.Range("E17").Value = '  -1.26%  '
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("D19").Value = '''4.853'
$ws.Range("E19").Value = '  -2.28%  '
$ws.Range("D20").Value = '''10.51'
$ws.Range("E20").Value = '  -4.36%  '
$ws.Range("D21").Value = '''188.37'
$ws.Range("E21").Value = '  -1.24%  '
$ws.Range("D22").Value = '''6.193'
$ws.Range("E22").Value = '  -1.15%  '
$ws.Range("D23").Value = '''1.006'
$ws.Range("E23").Value = '  -0.50%  '
$ws.Range("D24").Value = '''147.96'
$ws.Range("E24").Value = '  -1.32%  '
$ws.Range("D25").Value = '''0.1259'
$ws.Range("E25").Value = '  -4.04%  '
$ws.Range("D26").Value = '''7.595'
$ws.Range("E26").Value = '  -4.08%  '
$ws.Range("D27").Value = '''15.82'
$ws.Range("E27").Value = '  +0.22%  '
$ws.Range("D28").Value = '''0.06439'
$ws.Range("E28").Value = '  +1.48%  '
$ws.Range("D29").Value = '''1.339'
$ws.Range("E29").Value = '  -6.07%  '
$ws.Range("D30").Value = '''1.277'
$ws.Range("E30").Value = '  -3.94%  '
$ws.Range("D31").Value = '''3.522'
$ws.Range("D32").Value = '''3.483'
$ws.Range("E32").Value = '  -3.22%  '
$ws.Range("D33").Value = '''1.647'
$ws.Range("E33").Value = '  -2.03%  '
$ws.Range("D34").Value = '''1.004'
$ws.Range("E34").Value = '  -4.08%  '
$ws.Range("D35").Value = '''0.6073'
$ws.Range("E35").Value = '  -2.72%  '
$ws.Range("D36").Value = '''2.415'
$ws.Range("D37").Value = '''2.713'
$ws.Range("D38").Value = '''6.148'
$ws.Range("E38").Value = '  -1.54%  '
$ws.Range("D39").Value = '1.099.61'
$ws.Range("E39").Value = '  -2.16%  '
$ws.Range("D40").Value = '''0.01604'
$ws.Range("E40").Value = '  -2.40%  '
$ws.Range("D41").Value = '''0.8645'
$ws.Range("E41").Value = '  -2.28%  '
$ws.Range("E42").Value = '  -1.08%  '
$ws.Range("E43").Value = '  -0.74%  '
$ws.Range("D44").Value = '1.829.13'
$ws.Range("E44").Value = '  -0.78%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").Value = '''56.71'
$ws.Range("E45").Value = '  -1.50%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '''0.00000000107'
$ws.Range("E46").Value = '  -3.73%  '
$ws.Range("D47").Value = '''1.003'
$ws.Range("E47").Value = '  -0.45%  '
$ws.Range("D48").Value = '''0.05269'
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("D49").Value = '''7.955'
$ws.Range("E49").Value = '  -3.25%  '
$ws.Range("D50").Value = '''0.4271'
$ws.Range("E50").Value = '  -0.74%  '
$ws.Range("D51").Value = '''5.936'
$ws.Range("E51").Value = '  -2.51%  '
